$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Serial"
$ws.Range("K1").NumberFormat = "General"

$ws.Range("D2").Copy()
$ws.Range("K2:K18").PasteSpecial(-4122)

$ws.Range("K2").Value = 1
$ws.Range("K3").Formula = "=1+K2"
$ws.Range("K4:K10").FormulaR1C1 = "=1+R[-1]C"
